# This workbook is a weekly price log for "Feria Lagunitas de Puerto Montt -
# Mango". A new weekly record was inserted in the middle of the data table
# (before what used to be row 194), which pushes every subsequent record
# down by one row (old row 194 -> 195, ..., old row 235 -> 236).
#
# Net effect vs. the original file:
#   - One new row is inserted at row 194.
#   - Rows 194..235 (old) become rows 195..236 (new) with identical content.
#   - The new row 194 holds a brand-new record (Brasil / Primera / 44798).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 194; this shifts rows 194-235 down
# to 195-236 and also copies formatting (e.g. the date style on column D)
# down from the row above, same as Excel's native "Insert" command.
$ws.Rows(194).Insert()

# Populate the newly inserted row 194 with the new weekly record.
$ws.Range("A194").Value2 = 4
$ws.Range("B194").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C194").Value2 = "Los Lagos"
$ws.Range("D194").Value2 = 44798
$ws.Range("E194").Value2 = 10
$ws.Range("F194").Value2 = "Fruta"
$ws.Range("G194").Value2 = 100108
$ws.Range("H194").Value2 = "Tropicales y subtropicales"
$ws.Range("I194").Value2 = 100108002
$ws.Range("J194").Value2 = "Mango"
$ws.Range("K194").Value2 = "Sin especificar"
$ws.Range("L194").Value2 = "Primera"
$ws.Range("M194").Value2 = 120
$ws.Range("N194").Value2 = 13000
$ws.Range("O194").Value2 = 14000
$ws.Range("P194").Value2 = 13500
$ws.Range("Q194").Value2 = "$/bandeja 4 kilos"
$ws.Range("R194").Value2 = "Brasil"
$ws.Range("S194").Value2 = 3375
$ws.Range("T194").Value2 = 4
